$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "ddd\ dd/mm/yyyy"

# Row 8: continuation of suspend work, no additional effort recorded
$ws.Range("A8").Value = 41174
$ws.Range("A8").NumberFormat = $dateFormat
$ws.Range("B8").Value = 5.25
$ws.Range("D8").Value = "Continuation of implementation rtos.c. Implementation problems with first suspend operation; concept made but not proven or implemented yet"

# Row 9: first suspend is running
$ws.Range("A9").Value = 41176
$ws.Range("A9").NumberFormat = $dateFormat
$ws.Range("B9").Value = 1.5
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "First suspend is running. TC02: Idle is periodically interrupted by a single task and then continued"

$ws.Range("D9").Select()
